$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "Mean" column (BA) entirely; AZ ("Run 50") becomes the
#     last column and will be repurposed below as the new "Mean" column. ---
$ws.Range("BA1:BA14").Delete()

# --- Header: "Gen" -> "MaxFES" ---
$ws.Range("A1").Value = "MaxFES"

# --- Column A: generation counts -> MaxFES fractions (rows 3-14; row 2 stays 0) ---
$ws.Range("A3").Value  = 0.001
$ws.Range("A4").Value  = 0.01
$ws.Range("A5").Value  = 0.1
$ws.Range("A6").Value  = 0.2
$ws.Range("A7").Value  = 0.3
$ws.Range("A8").Value  = 0.4
$ws.Range("A9").Value  = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# --- Column AZ header: "Run 50" -> "Mean" (recomputed over Run 0..Run 49) ---
$ws.Range("AZ1").Value = "Mean"

# --- Column AZ values: recomputed mean of B:AY (Run 0 .. Run 49) per row ---
$ws.Range("AZ2").Value  = 15891024578.90219
$ws.Range("AZ3").Value  = 12153786732.82431
$ws.Range("AZ4").Value  = 1208406654.984627
$ws.Range("AZ5").Value  = 9517.010041990001
$ws.Range("AZ6").Value  = 3801.01988885
$ws.Range("AZ7").Value  = 3201.95020445
$ws.Range("AZ8").Value  = 3085.36051031
$ws.Range("AZ9").Value  = 3039.30371846
$ws.Range("AZ10").Value = 3004.09919664
$ws.Range("AZ11").Value = 2986.41496231
$ws.Range("AZ12").Value = 2985.08701203
$ws.Range("AZ13").Value = 2984.95499071
$ws.Range("AZ14").Value = 2984.79444316
